# Feature/redescente lrm hub (#112)
# Update the IG ValueSet metadata sheet:
#   - "Experimental" row: set the Value cell (B7) to the literal text "false"
#   - "Date" row: bump the Value cell (B8) to the new publication timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Writing the bare word "false" via .Value would be auto-coerced to a
# Boolean by Excel's input parser, but the source data models this column
# as plain text. Build the literal string via a text formula in a scratch
# cell, then paste-special as values-only so the destination cell ends up
# holding real text (not a boolean), keeping its existing cell style.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""false"""
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$scratch.ClearContents()

# Plain text update - no special coercion risk here.
$ws.Range("B8").Value = "2025-11-04T10:04:56+00:00"
